# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output counts.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    5  = 47
    8  = 2019
    11 = 4342
    13 = 280
    15 = 6
    19 = 67
    20 = 3137
    22 = 456
    24 = 16
    26 = 85
    28 = 7
    32 = 518
    33 = 1721
    34 = 263
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
